# Update the "as_of_utc" timestamp column (AA) on the data sheets to
# reflect the latest publish run.
$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-17 11:13:31"
$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $newTimestamp
    }
}
